$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.985.78"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.743.09"
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'250.56"
$ws.Range("E5").Value = "  +7.36%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5155"
$ws.Range("E7").Value = "  -2.04%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.2761"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.06197"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.740.80"
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07238"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'15.21"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.6517"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.637"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'77.91"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "26.012.58"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'11.87"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000006810"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "1.966.65"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.303"
$ws.Range("E22").Value = "  -0.42%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'8.688"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.396"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'136.11"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.512"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'15.29"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.788"
$ws.Range("E28").Value = "  -0.75%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'105.94"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'3.957"
$ws.Range("E30").Value = "  +5.13%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.08254"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'3.670"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04684"
$ws.Range("E33").Value = "  +3.72%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.656"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.004"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6262"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'2.732"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01617"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'1.925"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'0.9994"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "'100.52"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.7639"
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.3879"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.022"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.368"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1135"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.64"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.05229"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'30.81"
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.575"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'0.3446"
$ws.Range("E51").Value = "  -0.08%  "
